$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The author tried to move several course offerings from one period/time-slot
# to the next one (shifting the period-number label down a row) but, per the
# commit message, only partially succeeded: the period number moved down
# while the original course-listing text stayed attached to the wrong label
# in a few spots. Re-create that exact end-state cell by cell.

function Set-TextValue {
    param($Range, [string]$Text)
    # Force the cell to stay a text value even when it looks like a plain
    # number (e.g. "13"), matching the original workbook's string storage,
    # then restore the cell's original (default) style.
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

Set-TextValue $ws.Range("C4") "13"
$ws.Range("C5").Value = "14`nCSM-225-13810-Trottier-EAX-017`nCSM-199-13805-nan-EAX-015`nCSM-199-13805-nan-EAX-016"

Set-TextValue $ws.Range("B6") "9"
$ws.Range("B7").Value = "11`nCSM-305-09814-nan-EAX-015`nCSM-305-09814-nan-EAX-016"

Set-TextValue $ws.Range("E7") "11"
$ws.Range("E8").Value = "21`nCM-415-21803-Rounds-IDE-318`nCSM-199-11810-nan-EAX-015`nCSM-199-11810-nan-EAX-016"

$ws.Range("C8").Value = "21`nCM-415-21803-Rounds-IDE-318"
$ws.Range("C9").Value = "23`nCSM-300-21804-Michienzi-EAX-015`nCSM-300-21804-Michienzi-EAX-016"

$ws.Range("D8").Value = "20`nCM-400-20801-Robblee-IDE-317"
$ws.Range("D9").Value = "22`nCSM-220-20806-Michienzi-EAX-015`nCSM-220-20806-Michienzi-EAX-016"
